# The weekly logged price sheet gains a new record.
# A row is inserted at row 41 (pushing the former row 41 down to row 42,
# which keeps its original values), and row 41 is populated with the data
# that used to live in row 40. Row 40 itself is then updated with the
# latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; former row 41 shifts down to row 42
# and keeps its existing values/style untouched.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the data that previously sat in
# row 40 (same record, just relocated one row down).
$ws.Range("A41").Value2 = 1
$ws.Range("B41").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value2 = "Arica y Parinacota"
$ws.Range("D41").Value2 = 44362
$ws.Range("E41").Value2 = 15
$ws.Range("F41").Value2 = 100112052
$ws.Range("G41").Value2 = "Albahaca"
$ws.Range("H41").Value2 = "Sin especificar"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 250
$ws.Range("K41").Value2 = 2800
$ws.Range("L41").Value2 = 3000
$ws.Range("M41").Value2 = 2900
$ws.Range("N41").Value2 = "`$/paquete"
$ws.Range("O41").Value2 = "Región de Arica y Parinacota"
$ws.Range("P41").Value2 = 2900
$ws.Range("Q41").Value2 = 1
$ws.Range("R41").Value2 = "Hortaliza"

# Update row 40 with the latest week's observation.
$ws.Range("D40").Value2 = 44769
$ws.Range("J40").Value2 = 300
$ws.Range("K40").Value2 = 2500
$ws.Range("M40").Value2 = 2750
$ws.Range("P40").Value2 = 2750
